# Update the "Price" (column D) and a couple of "Volume(1h)" (column E)
# values on Sheet1 to reflect the refreshed symbol list, per the
# "Updated symbol list ... with GitHub Actions" commit.
#
# Prices are stored as text (to preserve an exact fixed-precision
# representation, e.g. trailing zeros like "0.001000"), so each numeric
# value is entered with a leading apostrophe - exactly like typing
# '241.48 into Excel - which keeps the cell as Text instead of letting
# Excel auto-convert it to a Number (which would silently drop
# significant trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'241.48"
$ws.Range("D3").Value  = "'21.79"
$ws.Range("D4").Value  = "'5.363"
$ws.Range("D5").Value  = "'0.05672"
$ws.Range("D6").Value  = "'3.408"
$ws.Range("D7").Value  = "'6.286"
$ws.Range("D8").Value  = "'0.8070"
$ws.Range("D9").Value  = "'0.8638"
$ws.Range("D10").Value = "'0.1432"
$ws.Range("D11").Value = "'0.07277"
$ws.Range("D12").Value = "'0.03024"
$ws.Range("D13").Value = "'0.03143"
$ws.Range("D15").Value = "'3.907"
$ws.Range("D16").Value = "'0.001587"
$ws.Range("D17").Value = "'0.04825"
$ws.Range("D18").Value = "'0.0005820"
$ws.Range("D19").Value = "'0.006347"
$ws.Range("D20").Value = "'0.001000"
$ws.Range("D22").Value = "'0.0001502"
$ws.Range("D23").Value = "'3.740"
$ws.Range("D24").Value = "'2.147"
$ws.Range("D27").Value = "'0.0004006"
$ws.Range("D40").Value = "'0.03787"

$ws.Range("D41").Value = "'0.006681"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("D42").Value = "'0.1047"
$ws.Range("D43").Value = "'0.002685"
$ws.Range("D44").Value = "'0.006832"
$ws.Range("D45").Value = "'0.00005619"

$ws.Range("D47").Value = "'0.5810"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"

$ws.Range("D49").Value = "'0.00002104"
$ws.Range("D50").Value = "'0.01012"
